# Book and media report content update Tue 12/31/2024
$wb = $excel.ActiveWorkbook

# --- Update the "Data" sheet: 2024 row (row 13) Book Count / Total Hours ---
$wsData = $wb.Worksheets.Item("Data")
$wsData.Range("B13").Value = 52
$wsData.Range("C13").Value = 632

# Update the stored selection on the Data sheet to D10
$wsData.Range("D10").Select()

# --- Update the "Graph" sheet selection to T56 and make it the active sheet ---
$wsGraph = $wb.Worksheets.Item("Graph")
$wsGraph.Activate()
$wsGraph.Range("T56").Select()
